# Threat Alert Report update - 2026-01-18 01:00
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Nile Air NP-231 ---
$ws.Range("C2").Value = "Nile Air NP-231"
$ws.Range("D2").Value = 15532
$ws.Range("F2").Value = -3313
$ws.Range("G2").Value = 30
$ws.Range("I2").Value = 0

# --- Row 3: flynas XY-895 (was MEDIUM THREAT, now LOW THREAT) ---
$ws.Range("C3").Value = "flynas XY-895"
$ws.Range("D3").Value = 16716
$ws.Range("F3").Value = -2129
$ws.Range("G3").Value = 40
$ws.Range("I3").Value = -10
$ws.Range("J3").Value = "LOW THREAT"
# Re-use the existing "LOW THREAT" (green) formatting from J2 instead of
# the old "MEDIUM THREAT" (yellow) formatting.
$ws.Range("J2").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = "LOW THREAT"

# --- Row 4: flyadeal F3-912 ---
$ws.Range("C4").Value = "flyadeal F3-912"
$ws.Range("D4").Value = 17623
$ws.Range("F4").Value = -1222

# --- Row 5 (old "flynas XY-895" row) is gone; data moved up to row 3 above ---
$ws.Rows.Item(5).Delete()

# --- Column J is narrower now ---
$ws.Columns.Item(10).ColumnWidth = 11.166666666666666

$excel.CutCopyMode = $false
